# Refresh cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price column (D) values that look like plain numbers are written with a
# leading apostrophe so Excel keeps them as text (preserving formats like
# trailing zeros / multi-dot thousands separators) instead of coercing them
# to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.106.70"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "1.625.41"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'214.17"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("E6").Value = "  +1.58%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "'20.33"
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("D11").Value = "'0.0847"
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").Value = "1.636.03"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "'64.68"
$ws.Range("E15").Value = "  -3.99%  "
$ws.Range("D16").Value = "27.105.16"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("D18").Value = "'217.17"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "'6.96"
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("E22").Value = "  -5.93%  "
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("D24").Value = "'148.10"
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -3.16%  "
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("D28").Value = "'15.63"
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("D33").Value = "1.348.88"
$ws.Range("E33").Value = "  +6.02%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("D36").Value = "'0.0177"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("D38").Value = "'0.858"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").Value = "'0.803"
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("E42").Value = "  +6.22%  "
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("D44").Value = "1.763.32"
$ws.Range("D45").Value = "'90.70"
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("E46").Value = "  +0.66%  "
$ws.Range("D47").Value = "'0.857"
$ws.Range("E47").Value = "  +29.01%  "
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("D49").Value = "'0.0514"
$ws.Range("E50").Value = "  +2.04%  "
$ws.Range("E51").Value = "  -0.35%  "
